$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header style from AC1 to the new header cells AD1:AF1
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill season record (Wins/Losses/Ties) for each player row
for ($r = 2; $r -le 48; $r++) {
    $ws.Cells.Item($r, 30).Value = 68
    $ws.Cells.Item($r, 31).Value = 93
    $ws.Cells.Item($r, 32).Value = 0
}
